# Swap the data in row 9 and row 10 (columns B:K), keeping column A (model name) untouched.
# Both rows are for "Mistral-7B-Instruct-v0.1"; row 9 was "text" and row 10 was "json_object".
# After the edit, row 9 should hold the "json_object" data and row 10 the "text" data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row9 = $ws.Range("B9:K9").Value2
$row10 = $ws.Range("B10:K10").Value2

$ws.Range("B9:K9").Value2 = $row10
$ws.Range("B10:K10").Value2 = $row9
